$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("M2").Value = 5.483544666666667
$ws.Range("N2").Value = 16.450634
$ws.Range("O2").Value = 0.4099269772514247
$ws.Range("P2").Value = 0.4099269772514247
$ws.Range("Q2").Value = 8.052625555660889
$ws.Range("R2").Value = 72.47363000094799
$ws.Range("S2").Value = 0.002098120825762494
$ws.Range("T2").Value = 0.002098120825762494
$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("O3").Value = 0.32690571008764
$ws.Range("P3").Value = 0.32690571008764
$ws.Range("Q3").Value = 6.421751730012666
$ws.Range("R3").Value = 57.79576557011399
$ws.Range("S3").Value = 0.001673194779700656
$ws.Range("T3").Value = 0.001673194779700656
$ws.Range("G4").Value = 1.468507333333333
$ws.Range("H4").Value = 4.405521999999999
$ws.Range("I4").Value = 0.005118279455112885
$ws.Range("J4").Value = 0.005118279455112885
$ws.Range("M4").Value = 2.750327
$ws.Range("N4").Value = 8.250980999999999
$ws.Range("O4").Value = 0.2056029999019452
$ws.Range("P4").Value = 0.2056029999019452
$ws.Range("Q4").Value = 4.038875368564666
$ws.Range("R4").Value = 36.34987831708199
$ws.Range("S4").Value = 0.001052333610307703
$ws.Range("T4").Value = 0.001052333610307703
$ws.Range("G5").Value = 1.468507333333333
$ws.Range("H5").Value = 4.405521999999999
$ws.Range("I5").Value = 0.005118279455112885
$ws.Range("J5").Value = 0.005118279455112885
$ws.Range("M5").Value = 0.770031
$ws.Range("N5").Value = 2.310093
$ws.Range("O5").Value = 0.05756431275899004
$ws.Range("P5").Value = 0.05756431275899004
$ws.Range("Q5").Value = 1.130796170394
$ws.Range("R5").Value = 10.177165533546
$ws.Range("S5").Value = 0.0002946302393420313
$ws.Range("T5").Value = 0.0002946302393420313
$ws.Range("I6").Value = 0.9046276674881553
$ws.Range("J6").Value = 0.9046276674881553
$ws.Range("M6").Value = 5.483544666666667
$ws.Range("N6").Value = 16.450634
$ws.Range("O6").Value = 0.4099269772514247
$ws.Range("P6").Value = 0.4099269772514247
$ws.Range("Q6").Value = 1423.257158476579
$ws.Range("R6").Value = 12809.31442628921
$ws.Range("S6").Value = 0.3708312852714264
$ws.Range("T6").Value = 0.3708312852714264
$ws.Range("I7").Value = 0.9046276674881553
$ws.Range("J7").Value = 0.9046276674881553
$ws.Range("O7").Value = 0.32690571008764
$ws.Range("P7").Value = 0.32690571008764
$ws.Range("S7").Value = 0.2957279500051408
$ws.Range("T7").Value = 0.2957279500051408
$ws.Range("I8").Value = 0.9046276674881553
$ws.Range("J8").Value = 0.9046276674881553
$ws.Range("M8").Value = 2.750327
$ws.Range("N8").Value = 8.250980999999999
$ws.Range("O8").Value = 0.2056029999019452
$ws.Range("P8").Value = 0.2056029999019452
$ws.Range("Q8").Value = 713.8489478705955
$ws.Range("R8").Value = 6424.640530835358
$ws.Range("S8").Value = 0.1859941622298641
$ws.Range("T8").Value = 0.1859941622298641
$ws.Range("I9").Value = 0.9046276674881553
$ws.Range("J9").Value = 0.9046276674881553
$ws.Range("M9").Value = 0.770031
$ws.Range("N9").Value = 2.310093
$ws.Range("O9").Value = 0.05756431275899004
$ws.Range("P9").Value = 0.05756431275899004
$ws.Range("Q9").Value = 199.861987021086
$ws.Range("R9").Value = 1798.757883189774
$ws.Range("S9").Value = 0.05207426998172382
$ws.Range("T9").Value = 0.05207426998172382
$ws.Range("G10").Value = 0.5890733333333333
$ws.Range("H10").Value = 1.76722
$ws.Range("I10").Value = 0.002053133730501083
$ws.Range("J10").Value = 0.002053133730501083
$ws.Range("M10").Value = 5.483544666666667
$ws.Range("N10").Value = 16.450634
$ws.Range("O10").Value = 0.4099269772514247
$ws.Range("P10").Value = 0.4099269772514247
$ws.Range("Q10").Value = 3.230209935275556
$ws.Range("R10").Value = 29.07188941748
$ws.Range("S10").Value = 0.0008416349040372502
$ws.Range("T10").Value = 0.0008416349040372502
$ws.Range("G11").Value = 0.5890733333333333
$ws.Range("H11").Value = 1.76722
$ws.Range("I11").Value = 0.002053133730501083
$ws.Range("J11").Value = 0.002053133730501083
$ws.Range("O11").Value = 0.32690571008764
$ws.Range("P11").Value = 0.32690571008764
$ws.Range("Q11").Value = 2.576005316126667
$ws.Range("R11").Value = 23.18404784514
$ws.Range("S11").Value = 0.0006711811400743419
$ws.Range("T11").Value = 0.0006711811400743419
$ws.Range("G12").Value = 0.5890733333333333
$ws.Range("H12").Value = 1.76722
$ws.Range("I12").Value = 0.002053133730501083
$ws.Range("J12").Value = 0.002053133730501083
$ws.Range("M12").Value = 2.750327
$ws.Range("N12").Value = 8.250980999999999
$ws.Range("O12").Value = 0.2056029999019452
$ws.Range("P12").Value = 0.2056029999019452
$ws.Range("Q12").Value = 1.620144293646667
$ws.Range("R12").Value = 14.58129864282
$ws.Range("S12").Value = 0.0004221304541908947
$ws.Range("T12").Value = 0.0004221304541908947
$ws.Range("G13").Value = 0.5890733333333333
$ws.Range("H13").Value = 1.76722
$ws.Range("I13").Value = 0.002053133730501083
$ws.Range("J13").Value = 0.002053133730501083
$ws.Range("M13").Value = 0.770031
$ws.Range("N13").Value = 2.310093
$ws.Range("O13").Value = 0.05756431275899004
$ws.Range("P13").Value = 0.05756431275899004
$ws.Range("Q13").Value = 0.45360472794
$ws.Range("R13").Value = 4.082442551460001
$ws.Range("S13").Value = 0.0001181872321985963
$ws.Range("T13").Value = 0.0001181872321985963
$ws.Range("G14").Value = 25.306101
$ws.Range("H14").Value = 75.91830299999999
$ws.Range("I14").Value = 0.0882009193262308
$ws.Range("J14").Value = 0.0882009193262308
$ws.Range("M14").Value = 5.483544666666667
$ws.Range("N14").Value = 16.450634
$ws.Range("O14").Value = 0.4099269772514247
$ws.Range("P14").Value = 0.4099269772514247
$ws.Range("Q14").Value = 138.767135172678
$ws.Range("R14").Value = 1248.904216554102
$ws.Range("S14").Value = 0.03615593625019856
$ws.Range("T14").Value = 0.03615593625019856
$ws.Range("G15").Value = 25.306101
$ws.Range("H15").Value = 75.91830299999999
$ws.Range("I15").Value = 0.0882009193262308
$ws.Range("J15").Value = 0.0882009193262308
$ws.Range("O15").Value = 0.32690571008764
$ws.Range("P15").Value = 0.32690571008764
$ws.Range("Q15").Value = 110.663048244879
$ws.Range("R15").Value = 995.9674342039109
$ws.Range("S15").Value = 0.02883338416272413
$ws.Range("T15").Value = 0.02883338416272413
$ws.Range("G16").Value = 25.306101
$ws.Range("H16").Value = 75.91830299999999
$ws.Range("I16").Value = 0.0882009193262308
$ws.Range("J16").Value = 0.0882009193262308
$ws.Range("M16").Value = 2.750327
$ws.Range("N16").Value = 8.250980999999999
$ws.Range("O16").Value = 0.2056029999019452
$ws.Range("P16").Value = 0.2056029999019452
$ws.Range("Q16").Value = 69.60005284502699
$ws.Range("R16").Value = 626.4004756052429
$ws.Range("S16").Value = 0.01813437360758251
$ws.Range("T16").Value = 0.01813437360758251
$ws.Range("G17").Value = 25.306101
$ws.Range("H17").Value = 75.91830299999999
$ws.Range("I17").Value = 0.0882009193262308
$ws.Range("J17").Value = 0.0882009193262308
$ws.Range("M17").Value = 0.770031
$ws.Range("N17").Value = 2.310093
$ws.Range("O17").Value = 0.05756431275899004
$ws.Range("P17").Value = 0.05756431275899004
$ws.Range("Q17").Value = 19.486482259131
$ws.Range("R17").Value = 175.378340332179
$ws.Range("S17").Value = 0.005077225305725599
$ws.Range("T17").Value = 0.005077225305725599
